$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the "IBKS스팩24호" IPO entry right after the
# existing "에이피알" row (row 2). This pushes every subsequent row
# down by one.
$ws.Range("A3").EntireRow.Insert()

$ws.Range("A3").Value = "IBKS스팩24호"
$ws.Range("B3").Value = "2024.01.17~01.18"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 2147483647
$ws.Range("F3").Value = "아이비케이투자증권"

# The table keeps a fixed window of rows, so the last existing row
# (old row 21, "와이바이오로직스", now shifted to row 22) drops off.
$ws.Range("A22").EntireRow.Delete()
